# Add multiple buy and sell swaps and add liquidity test
#
# Target changes (per commit diff):
#  - CalculationAmountAndPrice (sheet1): add a "Price" column (N) with a
#    per-row price formula, change rows 5 & 6 so they model additional
#    BUY/SELL swaps (A/B running-balance formulas), and move the cell
#    selection.
#  - CalculationLiquidity (sheet2) becomes the active tab/sheet, with a new
#    selected cell.
#  - Staking (sheet3) selected cell moves too.

$wb = $excel.ActiveWorkbook

$wsPrice = $wb.Worksheets.Item("CalculationAmountAndPrice")
$wsLiquidity = $wb.Worksheets.Item("CalculationLiquidity")
$wsStaking = $wb.Worksheets.Item("Staking")

# --- CalculationAmountAndPrice: new "Price" column ------------------------
$wsPrice.Range("N1").Value = "Price"

$wsPrice.Range("N2").Formula = "= M2/E2"
$wsPrice.Range("N3").Formula = "= M3/E3"
$wsPrice.Range("N4").Formula = "=(M4/E4) /1"
$wsPrice.Range("N5").Formula = "=(M5/E5) /1"
$wsPrice.Range("N6").Formula = "=M6/E6"

# --- CalculationAmountAndPrice: rework rows 5 & 6 (more buy/sell swaps) ---
$wsPrice.Range("A5").Formula = "=A4 -M4"
$wsPrice.Range("B5").Formula = "=B4 +M4+H4"

$wsPrice.Range("A6").Formula = "=A5 -M5"
$wsPrice.Range("B6").Formula = "=B5 + L5+H5"

# --- selection / active-sheet bookkeeping ----------------------------------
$wsPrice.Range("D20").Select()
$wsStaking.Range("C28").Select()

$wsLiquidity.Activate()
$wsLiquidity.Range("H21").Select()
